# Update TPM-derived values for Wnt1-Fzd8 ligand-receptor pair sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.027123
$ws.Range("H2").Value = 0.081369
$ws.Range("I2").Value = 0.0960827240265261
$ws.Range("J2").Value = 0.09608272402652611
$ws.Range("M2").Value = 1.824475333333333
$ws.Range("N2").Value = 5.473426
$ws.Range("O2").Value = 0.1906606574278047
$ws.Range("P2").Value = 0.2015451970524477
$ws.Range("Q2").Value = 0.04948524446599999
$ws.Range("R2").Value = 0.445367200194
$ws.Range("S2").Value = 0.01831919533035179
$ws.Range("T2").Value = 0.01936501154726215

# Row 3
$ws.Range("G3").Value = 0.027123
$ws.Range("H3").Value = 0.081369
$ws.Range("I3").Value = 0.0960827240265261
$ws.Range("J3").Value = 0.09608272402652611
$ws.Range("O3").Value = 0.6423822165107047
$ws.Range("P3").Value = 0.6790548829333741
$ws.Range("Q3").Value = 0.166727847546
$ws.Range("R3").Value = 1.500550627914
$ws.Range("S3").Value = 0.06172183322854617
$ws.Range("T3").Value = 0.06524544291575238

# Row 4
$ws.Range("G4").Value = 0.027123
$ws.Range("H4").Value = 0.081369
$ws.Range("I4").Value = 0.0960827240265261
$ws.Range("J4").Value = 0.09608272402652611
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.009389666666666666
$ws.Range("N4").Value = 0.028169
$ws.Range("O4").Value = 0.0009812355294625031
$ws.Range("P4").Value = 0.001037252838673693
$ws.Range("Q4").Value = 0.000254675929
$ws.Range("R4").Value = 0.002292083361
$ws.Range("S4").Value = 0.00009427978258236791
$ws.Range("T4").Value = 0.00009966207824401529

# Row 5
$ws.Range("G5").Value = 0.027123
$ws.Range("H5").Value = 0.081369
$ws.Range("I5").Value = 0.0960827240265261
$ws.Range("J5").Value = 0.09608272402652611
$ws.Range("M5").Value = 1.5503715
$ws.Range("N5").Value = 3.100743
$ws.Range("O5").Value = 0.1620163583726162
$ws.Range("P5").Value = 0.1141770910840848
$ws.Range("Q5").Value = 0.0420507261945
$ws.Range("R5").Value = 0.252304357167
$ws.Range("S5").Value = 0.01556697304929883
$ws.Range("T5").Value = 0.01097044593278365

# Row 6
$ws.Range("G6").Value = 0.027123
$ws.Range("H6").Value = 0.081369
$ws.Range("I6").Value = 0.0960827240265261
$ws.Range("J6").Value = 0.09608272402652611
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.03788966666666667
$ws.Range("N6").Value = 0.113669
$ws.Range("O6").Value = 0.003959532159411881
$ws.Range("P6").Value = 0.004185576091419648
$ws.Range("Q6").Value = 0.001027681429
$ws.Range("R6").Value = 0.009249132861000001
$ws.Range("S6").Value = 0.0003804426357469267
$ws.Range("T6").Value = 0.0004021615524838998

# Row 7
$ws.Range("I7").Value = 0.9039172759734738
$ws.Range("J7").Value = 0.9039172759734738
$ws.Range("M7").Value = 1.824475333333333
$ws.Range("N7").Value = 5.473426
$ws.Range("O7").Value = 0.1906606574278047
$ws.Range("P7").Value = 0.2015451970524477
$ws.Range("Q7").Value = 0.4655422484299999
$ws.Range("R7").Value = 4.18988023587
$ws.Range("S7").Value = 0.1723414620974529
$ws.Range("T7").Value = 0.1821801855051855

# Row 8
$ws.Range("I8").Value = 0.9039172759734738
$ws.Range("J8").Value = 0.9039172759734738
$ws.Range("O8").Value = 0.6423822165107047
$ws.Range("P8").Value = 0.6790548829333741
$ws.Range("S8").Value = 0.5806603832821585
$ws.Range("T8").Value = 0.6138094400176217

# Row 9
$ws.Range("I9").Value = 0.9039172759734738
$ws.Range("J9").Value = 0.9039172759734738
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.009389666666666666
$ws.Range("N9").Value = 0.028169
$ws.Range("O9").Value = 0.0009812355294625031
$ws.Range("P9").Value = 0.001037252838673693
$ws.Range("Q9").Value = 0.002395914295
$ws.Range("R9").Value = 0.021563228655
$ws.Range("S9").Value = 0.0008869557468801351
$ws.Range("T9").Value = 0.0009375907604296781

# Row 10
$ws.Range("I10").Value = 0.9039172759734738
$ws.Range("J10").Value = 0.9039172759734738
$ws.Range("M10").Value = 1.5503715
$ws.Range("N10").Value = 3.100743
$ws.Range("O10").Value = 0.1620163583726162
$ws.Range("P10").Value = 0.1141770910840848
$ws.Range("Q10").Value = 0.3956005437975
$ws.Range("R10").Value = 2.373603262785
$ws.Range("S10").Value = 0.1464493853233174
$ws.Range("T10").Value = 0.1032066451513011

# Row 11
$ws.Range("I11").Value = 0.9039172759734738
$ws.Range("J11").Value = 0.9039172759734738
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.03788966666666667
$ws.Range("N11").Value = 0.113669
$ws.Range("O11").Value = 0.003959532159411881
$ws.Range("P11").Value = 0.004185576091419648
$ws.Range("Q11").Value = 0.009668116794999999
$ws.Range("R11").Value = 0.08701305115499999
$ws.Range("S11").Value = 0.003579089523664954
$ws.Range("T11").Value = 0.003783414538935748

Write-Output "Updated Wnt1-Fzd8 TPM values"
